$wb = $excel.ActiveWorkbook

# Sheet ALC row 20
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Sheet ALC row 35
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

# Sheet ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8281.137000000001
$ws.Range("I62").Value = 10652.23
$ws.Range("J62").Value = 4856.222
$ws.Range("K62").Value = 10652.23
$ws.Range("L62").Value = 4856.222
$ws.Range("M62").Value = -10028.23
$ws.Range("N62").Value = -6104.222

# Sheet ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 8281.137000000001
$ws.Range("I65").Value = 10652.23
$ws.Range("J65").Value = 4856.222
$ws.Range("K65").Value = 53261.14999999999
$ws.Range("L65").Value = 24281.11
$ws.Range("M65").Value = -50141.14999999999
$ws.Range("N65").Value = -30521.11

# Sheet ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6363.8667
$ws.Range("I76").Value = 6929
$ws.Range("J76").Value = 4934.4116
$ws.Range("K76").Value = 6929
$ws.Range("L76").Value = 4934.4116
$ws.Range("M76").Value = -6614
$ws.Range("N76").Value = -5564.4116

# Sheet ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 6363.8667
$ws.Range("I79").Value = 6929
$ws.Range("J79").Value = 4934.4116
$ws.Range("K79").Value = 6929
$ws.Range("L79").Value = 4934.4116
$ws.Range("M79").Value = -5837
$ws.Range("N79").Value = -7118.4116

# Sheet ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2949.1345
$ws.Range("I61").Value = 1864.7931
$ws.Range("J61").Value = 4316.3477
$ws.Range("K61").Value = 1864.7931
$ws.Range("L61").Value = 4316.3477
$ws.Range("M61").Value = -1652.7931
$ws.Range("N61").Value = -4740.3477

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2973.831
$ws.Range("I132").Value = 2693.8086
$ws.Range("J132").Value = 3522.2083
$ws.Range("K132").Value = 8081.425799999999
$ws.Range("L132").Value = 10566.6249
$ws.Range("M132").Value = -5551.425799999999
$ws.Range("N132").Value = -15626.6249

# Sheet ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2949.1345
$ws.Range("I136").Value = 1864.7931
$ws.Range("J136").Value = 4316.3477
$ws.Range("K136").Value = 5594.379300000001
$ws.Range("L136").Value = 12949.0431
$ws.Range("M136").Value = -3044.379300000001
$ws.Range("N136").Value = -18049.0431

# Sheet BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1161.0571
$ws.Range("I94").Value = 785.6087
$ws.Range("J94").Value = 1880.6666
$ws.Range("K94").Value = 785.6087
$ws.Range("L94").Value = 1880.6666
$ws.Range("M94").Value = -334.6087
$ws.Range("N94").Value = -2782.6666

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 33323.89
$ws.Range("I134").Value = 43966.4
$ws.Range("J134").Value = 9136.362999999999
$ws.Range("K134").Value = 131899.2
$ws.Range("L134").Value = 27409.089
$ws.Range("M134").Value = -129364.2
$ws.Range("N134").Value = -32479.089

# Sheet CRP row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 24333.334
$ws.Range("J50").Value = 24333.334
$ws.Range("L50").Value = 24333.334
$ws.Range("N50").Value = -25583.334

# Sheet CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2605.543
$ws.Range("I132").Value = 1176.4375
$ws.Range("J132").Value = 3809
$ws.Range("K132").Value = 3529.3125
$ws.Range("L132").Value = 11427
$ws.Range("M132").Value = -999.3125
$ws.Range("N132").Value = -16487

# Sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1790.3243
$ws.Range("I134").Value = 1056.8636
$ws.Range("J134").Value = 2866.0667
$ws.Range("K134").Value = 3170.5908
$ws.Range("L134").Value = 8598.2001
$ws.Range("M134").Value = -635.5907999999999
$ws.Range("N134").Value = -13668.2001

# Sheet CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3080.476
$ws.Range("I3").Value = 2947.611
$ws.Range("J3").Value = 3877.6667
$ws.Range("K3").Value = 8842.832999999999
$ws.Range("L3").Value = 11633.0001
$ws.Range("M3").Value = -8730.832999999999
$ws.Range("N3").Value = -11857.0001

# Sheet CUL row 13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 122.5
$ws.Range("I13").Value = 122.5
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 367.5
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -199.5
$ws.Range("N13").ClearContents()

# Sheet CUL row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1933.9048
$ws.Range("I64").Value = 845.7778
$ws.Range("J64").Value = 2750
$ws.Range("K64").Value = 2537.3334
$ws.Range("L64").Value = 8250
$ws.Range("M64").Value = -2267.3334
$ws.Range("N64").Value = -8790

# Sheet CUL row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 1933.9048
$ws.Range("I67").Value = 845.7778
$ws.Range("J67").Value = 2750
$ws.Range("K67").Value = 2537.3334
$ws.Range("L67").Value = 8250
$ws.Range("M67").Value = -1601.3334
$ws.Range("N67").Value = -10122

# Sheet CUL row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 576.8
$ws.Range("I114").Value = 445.77777
$ws.Range("J114").Value = 773.3333
$ws.Range("K114").Value = 1337.33331
$ws.Range("L114").Value = 2319.9999
$ws.Range("M114").Value = 1916.66669
$ws.Range("N114").Value = -8827.999899999999

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 845.4143
$ws.Range("I131").Value = 410.33334
$ws.Range("J131").Value = 909.60657
$ws.Range("K131").Value = 1231.00002
$ws.Range("L131").Value = 2728.81971
$ws.Range("M131").Value = 3808.99998
$ws.Range("N131").Value = -12808.81971

# Sheet CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 1931.9445
$ws.Range("I133").Value = 2213.182
$ws.Range("J133").Value = 1808.2
$ws.Range("K133").Value = 6639.545999999999
$ws.Range("L133").Value = 5424.6
$ws.Range("M133").Value = -1579.545999999999
$ws.Range("N133").Value = -15544.6

# Sheet CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 3468.0667
$ws.Range("I136").Value = 2249
$ws.Range("J136").Value = 5906.2
$ws.Range("K136").Value = 6747
$ws.Range("L136").Value = 17718.6
$ws.Range("M136").Value = -1647
$ws.Range("N136").Value = -27918.6

# Sheet CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2592.84
$ws.Range("I137").Value = 2026.3636
$ws.Range("J137").Value = 3037.9285
$ws.Range("K137").Value = 6079.0908
$ws.Range("L137").Value = 9113.7855
$ws.Range("M137").Value = -979.0907999999999
$ws.Range("N137").Value = -19313.7855

# Sheet GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3542.7144
$ws.Range("I126").Value = 3466.5
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 10399.5
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -7929.5
$ws.Range("N126").Value = -16940

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4524.5312
$ws.Range("I132").Value = 6506.1665
$ws.Range("J132").Value = 3335.55
$ws.Range("K132").Value = 19518.4995
$ws.Range("L132").Value = 10006.65
$ws.Range("M132").Value = -16988.4995
$ws.Range("N132").Value = -15066.65

# Sheet GSM row 139
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 24228
$ws.Range("J139").Value = 24228
$ws.Range("L139").Value = 24228
$ws.Range("N139").Value = -34508

# Sheet LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3827.7778
$ws.Range("I40").Value = 3916.6667
$ws.Range("J40").Value = 3650
$ws.Range("K40").Value = 3916.6667
$ws.Range("L40").Value = 3650
$ws.Range("M40").Value = -3780.6667
$ws.Range("N40").Value = -3922

# Sheet WVR row 58
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 12475
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 12475
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 12475
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -13091

# Sheet WVR row 101
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 12648.333
$ws.Range("J101").Value = 12648.333
$ws.Range("L101").Value = 12648.333
$ws.Range("N101").Value = -19138.333

# Sheet WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 339.96
$ws.Range("I113").Value = 270.89474
$ws.Range("J113").Value = 558.6667
$ws.Range("K113").Value = 812.6842200000001
$ws.Range("L113").Value = 1676.0001
$ws.Range("M113").Value = 1357.31578
$ws.Range("N113").Value = -6016.0001

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 20568.389
$ws.Range("I132").Value = 30603.559
$ws.Range("J132").Value = 3508.6
$ws.Range("K132").Value = 91810.677
$ws.Range("L132").Value = 10525.8
$ws.Range("M132").Value = -89280.677
$ws.Range("N132").Value = -15585.8

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 31252646
$ws.Range("I136").Value = 62502116
$ws.Range("J136").Value = 3175.625
$ws.Range("K136").Value = 187506348
$ws.Range("L136").Value = 9526.875
$ws.Range("M136").Value = -187503798
$ws.Range("N136").Value = -14626.875

